$wb = $excel.ActiveWorkbook

# Sheet 1: "HSV Log OLS"
$ws1 = $wb.Worksheets.Item("HSV Log OLS")
$ws1.Range("A2").Value = 0.011
$ws1.Range("B2").Value = 1.04
$ws1.Range("C2").Value = 0.913

# Sheet 2: "HSV PPML"
$ws2 = $wb.Worksheets.Item("HSV PPML")
$ws2.Range("A2").Value = -0.01
$ws2.Range("B2").Value = 0.813
$ws2.Range("C2").Value = 0.915

# Sheet 3: "HSVT NLLSQ"
$ws3 = $wb.Worksheets.Item("HSVT NLLSQ")
$ws3.Range("A2").Value = -0.019
$ws3.Range("B2").Value = 0.72
$ws3.Range("C2").Value = 1849.61
$ws3.Range("D2").Value = 0.014
